# Auto-generated edit script for sample_deployment.xlsx (Vscs sheet)
# Adds a 'Router ID' field (row 18) for IPv6 VSC support, shifting subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# --- Step 1: remove existing comments so they don't end up anchored to the wrong row
#            once row 18 is inserted and everything below it shifts down. ---
$existingCommentRefs = @(
    "A5", "A6", "A7", "A8", "A9", "A11", "A12", "A13", "A14", "A15", "A17", "A18", "A19", 
    "A20", "A22", "A23", "A25", "A26", "A27", "A28", "A29", "A30", "A31", "A33", "A34", 
    "A35", "A36", "A37", "A38", "A40", "A41", "A42", "A43", "A44", "A45", "A46", "A47", 
    "A48", "A49", "A50", "A51", "A52", "A54", "A55", "A56", "A57", "A58", "A60", "A61", 
    "A62", "A63", "A64", "A65", "A67", "A68", "A69" 
)
foreach ($ref in $existingCommentRefs) {
    $cell = $ws.Range($ref)
    if ($cell.Comment) {
        $cell.Comment.Delete()
    }
}

# --- Step 2: insert the new row 18 ('Router ID'); Excel shifts rows 18-69 down to 19-70,
#            carrying along merged cells and data validations automatically. ---
$ws.Rows("18:18").Insert()

# Row 18 inherits whole-row formatting from the Insert, but column B/C should keep the
# non-header (value-cell) style used throughout the sheet - copy it down from row 19
# (formerly row 18, 'VM name') which already has the correct style.
$ws.Range("B19:C19").Copy($ws.Range("B18:C18"))

$ws.Range("A18").Value = "Router ID"

# --- Step 3: re-create all comments at their (possibly shifted) positions. ---
$newComments = @{
    'A5' = 'Hostname of the VSC instance'
    'A6' = 'The BOF/Management IP address of the controller'
    'A7' = 'Management network prefix length'
    'A8' = 'Gateway IP on the Management network'
    'A9' = 'Static Route list to be configured on the management/BOF interface. Define as empty list if no static routes are to be configured. [default: [ 0.0.0.0/1, 128.0.0.0/1 ]] (List items separated by comma.)'
    'A11' = 'The Control/Data IP address of the controller'
    'A12' = 'Control network prefix length'
    'A13' = 'The VLAN ID for the uplink connection of the controller. This field is deprecated and will be removed in future releases. Use ctrl_ip_vprn_list instead. [default: 0]'
    'A14' = 'The Control/Data VPRN IP address of the controller. This field is deprecated and will be removed in future releases. Use ctrl_ip_vprn_list instead.'
    'A15' = 'Comma Seperated list of control IP VPRNs. Each item in list is of format VLAN_ID:Control_IP_ADDR/NETMASK_LENGTH . Both Control IP address and Netmask length are optional (List items separated by comma.)'
    'A17' = 'Required for BGP pairing with peers [default: ]'
    'A18' = 'Required IPv4 address when using an IPv6 system IP address [default: (system_ip)]'
    'A19' = 'Name of the virtual machine on the Hypervisor/vCenter Server. [default: (Hostname)]'
    'A20' = 'Name of the system if other than hostname [default: (Hostname)]'
    'A21' = 'Unique username used to identify this VSC in its XMPP connection with VSD [default: vsc1]'
    'A23' = 'Type of hypervisor environment where VMs will be instantiated. Use ''none'' when skipping predeploy.'
    'A24' = 'Hostname or IP address of the hypervisor where VM  will be instantiated. In the case of deployment in a vCenter environment, this will be the FQDN of the vCenter Server'
    'A26' = 'Network Bridge used for the management interface of a component or the BOF interface on VSC. This will be a Distributed Virtual PortGroup (DVPG) when deploying on vCenter or a Linux network bridge when deploying on KVM. This field can be overridden by defining the management network bridge separately in the component configuration. Defaults to the global setting [default: (global Bridge interface)]'
    'A27' = 'Network Bridge used for the data path of a component or the Control interface on VSC. This will be a Distributed Virtual PortGroup (DVPG) when deploying on vCenter or a Linux network bridge when deploying on KVM. [default: (global Bridge interface)]'
    'A28' = 'FQDN of the VSD or VSD cluster for this VSC'
    'A29' = 'Private Management IP Address of VSC instances'
    'A30' = 'Private Control IP Address of VSC Instances'
    'A31' = 'Private Data Gateway IP Address of VSC Instances'
    'A32' = 'List of route reflector IP addresses if present (List items separated by comma.)'
    'A34' = 'Name of the vCenter Datacenter on which the VSC VM will be deployed. Defaults to the common vCenter Datacenter Name if not defined here. [default: (global vCenter Datacenter Name)]'
    'A35' = 'Name of the vCenter Cluster on which the VSC VM will be deployed. Defaults to the common vCenter Cluster Name if not defined here. [default: (global vCenter Cluster Name)]'
    'A36' = 'Requires ovftool 4.3. Reference to the host on the vCenter cluster on which to deploy Nuage components [default: (global vCenter Host Reference)]'
    'A37' = 'Name of the vCenter Datastore on which the VSC VM will be deployed. Defaults to the common vCenter Datastore Name if not defined here. [default: (global vCenter Datastore Name)]'
    'A38' = 'Optional path to a folder defined on vCenter where VM will be instantiated [default: (global vCenter VM folder)]'
    'A39' = 'Optional path to a hosts and clusters folder defined on vCenter where VM will be instantiated'
    'A41' = 'Name of image installed on OpenStack for VSC'
    'A42' = 'Name of instance flavor installed on OpenStack for VSC'
    'A43' = 'Name of availability zone on OpenStack for VSC'
    'A44' = 'Name of management network on OpenStack for VSC'
    'A45' = 'Name of management subnet on OpenStack for VSC'
    'A46' = 'Name for Mgmt interface'
    'A47' = 'Set of security groups to associate with Mgmt interface (List items separated by comma.)'
    'A48' = 'Name of control network on OpenStack for VSC'
    'A49' = 'Name of control subnet on OpenStack for VSC'
    'A50' = 'Name for Control interface'
    'A51' = 'Set of security groups to associate with Control interface (List items separated by comma.)'
    'A52' = 'Name for Mgmt interface'
    'A53' = 'Set of security groups to associate with Mgmt interface (List items separated by comma.)'
    'A55' = 'Used in postdeploy and health workflows as expected values if non-zero [default: 0]'
    'A56' = 'Used in postdeploy and health workflows as expected values if non-zero [default: 0]'
    'A57' = 'Used in postdeploy and health workflows as expected values if non-zero [default: 0]'
    'A58' = 'Used in postdeploy and health workflows as expected values if non-zero [default: 0]'
    'A59' = 'Used in postdeploy and health workflows as expected values if non-zero [default: 0]'
    'A61' = 'Ejabberd user id used to create the certificate'
    'A62' = 'Path to VSC certificate key pem file'
    'A63' = 'Path to VSC certificate pem file'
    'A64' = 'Path to CA certificate pem file'
    'A65' = 'XMPP domain used in custom certificates'
    'A66' = 'Name of the credentials set for the vsc'
    'A68' = 'Cpuset information for cpu pinning on KVM. For example, VSC requires 4 cores and sample values will be of the form [ 0, 1, 2, 3 ] (List items separated by comma.)'
    'A69' = 'Enables hardening configuration on VSC [default: True]'
    'A70' = 'Paths to files that can be optionally applied for additional VSC configuration (List items separated by comma.)'
}
foreach ($ref in $newComments.Keys) {
    $ws.Range($ref).AddComment($newComments[$ref]) | Out-Null
}

